# Append a new data row (row 80) to the "NEW" sheet, mirroring the
# structure/format of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$row = 80

# Columns that hold numeric-looking text in the source data (Caso, F. De
# Reclamo, Direccion, Comuna, OT, Proveedor Asignado, Estado,
# Observaciones, Tipo de tarea, Equipo, Tipo de Elemento, Operacion,
# Zona) must be forced to Text so Excel does not silently reinterpret
# them as numbers/dates. NumberFormat is reset back to the default
# ("Normal" style) afterwards so the cell keeps plain/default formatting
# like every other data row, only the underlying value stays textual.
$textValues = @{
    1  = "-582"
    2  = "9/4/2025"
    3  = "Vilela 4019"
    4  = "12"
    5  = "809454353"
    6  = "NEW"
    7  = "Pendiente"
    8  = "Poste telefonico propio quebrado en base"
    10 = "Cambio"
    11 = "Sin equipos"
    12 = "Poste"
    15 = "Saavedra"
    16 = "Capital Norte"
}

foreach ($col in $textValues.Keys) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$col]
    $cell.Style = "Normal"
}

# Numeric columns keep their native numeric type.
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 13).Value = -58.485872
$ws.Cells.Item($row, 14).Value = -34.552645
